# Update the "Förändrad" date column (C) for rows 2-5 from 2023-10-08 to 2023-10-09
# (serial date 45207 -> 45208), preserving existing cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..5) {
    $cell = $ws.Cells.Item($row, 3)   # Column C
    $cell.Value = 45208
}
